$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(43, 1).Value = "Record"
$ws.Cells.Item(43, 2).Value = "Balanço Geral"
$ws.Cells.Item(43, 3).Value = "CCZ"
$ws.Cells.Item(43, 4).Value = "2025-04-04T11:57"
$ws.Cells.Item(43, 5).Value = "Positivo"
$ws.Cells.Item(43, 6).Value = "Cachorro e gato são resgatados em ação conjunta da Polícia Civil e CCZ. Repórter *ao vivo* do canil do CCZ. Entrevista com veterinário do CCZ, José Leonardo, que explicou como foi a ação. Exibido vídeo dos animais acorrentados. Esses não foram encontrados. Inicialmente, a denúncia era de maus tratos a animais. Eles estavam em ambiente insalubre. Ele também convidou as pessoas a aderir à adoção responsável. Maus tratos é crime. Abril Laranja é o mês de conscientização e combate aos maus-tratos a animais.  "

$ws.Cells.Item(44, 1).Value = "Record"
$ws.Cells.Item(44, 2).Value = "Balanço Geral"
$ws.Cells.Item(44, 3).Value = "Saúde"
$ws.Cells.Item(44, 4).Value = "2025-04-04T12:27"
$ws.Cells.Item(44, 5).Value = "Neutro"
$ws.Cells.Item(44, 6).Value = "Ambulância pega fogo em Campos. Repórter *ao vivo*. Veículo estava trafegando na Av. 28 de Março. Grande susto, mas ninguém se feriu. "
